# Apply the updated distractor-analysis numbers to the mv_person_all workbook.
# Sheets: OM, NV, NR, ND, ALL, summary  (ND is untouched)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# OM sheet: percentages change, and a new "4 missing" row is appended.
# ---------------------------------------------------------------------------
$om = $wb.Worksheets.Item("OM")
$om.Range("B2").Value = 80.19
$om.Range("B3").Value = 18.29
$om.Range("B4").Value = 1.33
$om.Range("B5").Value = 0.14
$om.Range("A6").Value = "'4"
$om.Range("B6").Value = 0.05

# ---------------------------------------------------------------------------
# NV sheet: percentages change only.
# ---------------------------------------------------------------------------
$nv = $wb.Worksheets.Item("NV")
$nv.Range("B2").Value = 94.38
$nv.Range("B3").Value = 5.52
$nv.Range("B4").Value = 0.1

# ---------------------------------------------------------------------------
# NR sheet: percentages change only.
# ---------------------------------------------------------------------------
$nr = $wb.Worksheets.Item("NR")
$nr.Range("B2").Value = 59.76
$nr.Range("B3").Value = 26.81
$nr.Range("B4").Value = 8.62
$nr.Range("B5").Value = 2.48
$nr.Range("B6").Value = 1.24
$nr.Range("B7").Value = 0.57
$nr.Range("B8").Value = 0.43
$nr.Range("B9").Value = 0.1

# ---------------------------------------------------------------------------
# ND sheet: unchanged.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# ALL sheet: percentages change, and the last "8 missing" row is removed.
# ---------------------------------------------------------------------------
$all = $wb.Worksheets.Item("ALL")
$all.Range("B2").Value = 43.29
$all.Range("B3").Value = 35.86
$all.Range("B4").Value = 13.48
$all.Range("B5").Value = 4.38
$all.Range("B7").Value = 0.62
$all.Range("B8").Value = 0.48
$all.Range("B9").Value = 0.1
$all.Rows.Item(10).Delete()

# ---------------------------------------------------------------------------
# summary sheet: distractor-analysis row 6 references the new max category
# (4 instead of 3 missing responses) and the new max count (7 instead of 8).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("summary")
$summary.Range("A6").Value = "'4"
$summary.Range("E6").Value = "'7"
